$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume/1h change (E) columns for the refreshed snapshot.
# Column D cells whose new value looks like a plain decimal number need the source
# column pre-formatted as Text so Excel keeps storing/display the literal digit string
# (e.g. "1.00", "20.32") instead of silently coercing it to a floating-point number.

$ws.Range("D2").Value = "58.116.48"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "3.137.99"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.58"
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.26"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E8").Value = "  +12.74%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.426"
$ws.Range("E10").Value = "  +7.58%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("D13").Value = "3.679.31"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.78"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("E15").Value = "  +5.44%  "
$ws.Range("D16").Value = "58.220.87"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.27"
$ws.Range("E17").Value = "  +7.23%  "
$ws.Range("D18").Value = "3.137.13"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  +4.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.21"
$ws.Range("E20").Value = "  +5.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.10"
$ws.Range("E21").Value = "  +8.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.46"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +11.55%  "
$ws.Range("D29").Value = "0.0₃0883"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("E30").Value = "  +6.30%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  +7.33%  "
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.76"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("E36").Value = "  +5.30%  "
$ws.Range("E37").Value = "  +11.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.65"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E39").Value = "  +5.78%  "
$ws.Range("D40").Value = "2.636.41"
$ws.Range("E40").Value = "  +10.47%  "
$ws.Range("E41").Value = "  +6.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0676"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.99"
$ws.Range("E43").Value = "  +6.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.701"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  +5.11%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").Value = "  +5.02%  "
$ws.Range("E48").Value = "  +11.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.976"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.32"
$ws.Range("E50").Value = "  +3.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.749"
$ws.Range("E51").Value = "  -0.68%  "
